$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.701.43"
$ws.Range("E2").Value = "'  +0.78%  "
$ws.Range("D3").Value = "'2.483.56"
$ws.Range("E3").Value = "'  +0.13%  "
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'587.06"
$ws.Range("E5").Value = "'  +0.26%  "
$ws.Range("D6").Value = "'174.55"
$ws.Range("E6").Value = "'  +0.99%  "
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "'  -0.49%  "
$ws.Range("D9").Value = "'0.144"
$ws.Range("E9").Value = "'  +4.26%  "
$ws.Range("E10").Value = "'  -1.44%  "
$ws.Range("D11").Value = "'4.94"
$ws.Range("E11").Value = "'  +0.23%  "
$ws.Range("D12").Value = "'0.332"
$ws.Range("E12").Value = "'  -0.41%  "
$ws.Range("D13").Value = "'2.937.82"
$ws.Range("E13").Value = "'  +0.20%  "
$ws.Range("D14").Value = "'25.22"
$ws.Range("E14").Value = "'  -1.40%  "
$ws.Range("D15").Value = "'67.773.70"
$ws.Range("E15").Value = "'  +1.02%  "
$ws.Range("D16").Value = "'0.0000169"
$ws.Range("E16").Value = "'  -0.70%  "
$ws.Range("D17").Value = "'2.493.20"
$ws.Range("E17").Value = "'  +0.38%  "
$ws.Range("D18").Value = "'10.77"
$ws.Range("E18").Value = "'  -1.77%  "
$ws.Range("D19").Value = "'7.37"
$ws.Range("E19").Value = "'  -2.74%  "
$ws.Range("D20").Value = "'346.33"
$ws.Range("E20").Value = "'  -1.22%  "
$ws.Range("D21").Value = "'4.09"
$ws.Range("E21").Value = "'  +1.48%  "
$ws.Range("E22").Value = "'  -0.03%  "
$ws.Range("D23").Value = "'70.67"
$ws.Range("E23").Value = "'  +2.36%  "
$ws.Range("D24").Value = "'4.17"
$ws.Range("E24").Value = "'  -1.58%  "
$ws.Range("D25").Value = "'1.67"
$ws.Range("E25").Value = "'  -7.49%  "
$ws.Range("D26").Value = "'8.82"
$ws.Range("E26").Value = "'  -4.17%  "
$ws.Range("D27").Value = "'2.611.31"
$ws.Range("E27").Value = "'  +0.23%  "
$ws.Range("E28").Value = "'  -0.21%  "
$ws.Range("D29").Value = "'0.0₃0889"
$ws.Range("E29").Value = "'  -2.47%  "
$ws.Range("D30").Value = "'498.05"
$ws.Range("E30").Value = "'  -1.88%  "
$ws.Range("D31").Value = "'7.69"
$ws.Range("E31").Value = "'  -0.20%  "
$ws.Range("E32").Value = "'  -0.51%  "
$ws.Range("D33").Value = "'1.75"
$ws.Range("E33").Value = "'  -1.00%  "
$ws.Range("E34").Value = "'  +0.00%  "
$ws.Range("D35").Value = "'164.56"
$ws.Range("E35").Value = "'  +1.25%  "
$ws.Range("E36").Value = "'  +1.20%  "
$ws.Range("E37").Value = "'  -0.38%  "
$ws.Range("D38").Value = "'18.24"
$ws.Range("E38").Value = "'  +0.60%  "
$ws.Range("E39").Value = "'  -0.04%  "
$ws.Range("E40").Value = "'  -3.04%  "
$ws.Range("D41").Value = "'1.71"
$ws.Range("E41").Value = "'  +1.43%  "
$ws.Range("D42").Value = "'0.323"
$ws.Range("E42").Value = "'  -1.88%  "
$ws.Range("D43").Value = "'4.75"
$ws.Range("E43").Value = "'  -1.55%  "
$ws.Range("D44").Value = "'2.37"
$ws.Range("E44").Value = "'  -1.03%  "
$ws.Range("D45").Value = "'147.93"
$ws.Range("E45").Value = "'  +3.21%  "
$ws.Range("D46").Value = "'3.51"
$ws.Range("E46").Value = "'  +0.72%  "
$ws.Range("D47").Value = "'0.509"
$ws.Range("E47").Value = "'  -1.14%  "
$ws.Range("D48").Value = "'0.0₆0252"
$ws.Range("E48").Value = "'  -4.38%  "
$ws.Range("E49").Value = "'  -0.21%  "
$ws.Range("D50").Value = "'1.55"
$ws.Range("E50").Value = "'  -1.37%  "
$ws.Range("E51").Value = "'  -1.26%  "
